$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns keep their original text (string) cell type,
# matching the source data which stores prices/percentages as text
# (e.g. "29.419.25" uses dots as thousands separators, not a real number).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.419.25'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.849.23'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '240.86'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '0.6326'
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.478.55'
$ws.Range("E8").Value = '  +88.13%  '
$ws.Range("D9").Value = '0.07582'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = '0.2971'
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").Value = '24.64'
$ws.Range("E11").Value = '  +1.29%  '
$ws.Range("D12").Value = '3.579.20'
$ws.Range("E12").Value = '  +71.42%  '
$ws.Range("D13").Value = '0.07711'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '4.985'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '0.6840'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '82.92'
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.000009966'
$ws.Range("E17").Value = '  +4.62%  '
$ws.Range("D18").Value = '6.165'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '29.440.46'
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").Value = '231.69'
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("D21").Value = '12.50'
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '7.576'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("D24").Value = '0.9998'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '155.13'
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("D26").Value = '0.1386'
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("D27").Value = '8.428'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("D29").Value = '1.467'
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").Value = '0.05793'
$ws.Range("E30").Value = '  -3.17%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").Value = '4.127'
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("B33").Value = 'RocketPoolETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D33").Value = '3.547.32'
$ws.Range("E33").Value = '  +77.02%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.021'
$ws.Range("E34").Value = '  -1.26%  '
$ws.Range("D35").Value = '1.869'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '1.156'
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("D37").Value = '0.7166'
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("D38").Value = '2.592'
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("D39").Value = '1.249.98'
$ws.Range("E39").Value = '  +3.96%  '
$ws.Range("D40").Value = '2.792'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = '0.01806'
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").Value = '0.8996'
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("D43").Value = '6.080'
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '101.67'
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D46").Value = '66.97'
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("D47").Value = '7.203'
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("D48").Value = '9.146'
$ws.Range("E48").Value = '  +0.50%  '
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("D51").Value = '0.1124'
$ws.Range("E51").Value = '  -0.39%  '
